$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "Hyperparameter Optimization": fill in previously-blank S:V columns
# (Average RMSE / Average MAE / Average MAPE / Average Pearson R2) for the
# rows that now have results.
# -------------------------------------------------------------------------
$wsHP = $wb.Worksheets.Item("Hyperparameter Optimization")

$hpData = @{
    53 = @(1.1521635149665499, 0.89770011856768095, 0.111066541482916, 0.96599362739920902)
    54 = @(1.92041390669594, 1.5638605476249701, 0.43751975447228802, 0.89923174406364903)
    55 = @(1.5953117484788699, 1.3134364638165901, 0.32637145109844501, 0.92623041087830604)
    57 = @(2.0477885452531499, 1.60133084206791, 0.10694157601608099, 0.94937875157241403)
    58 = @(2.7738462629876399, 2.1520019190546198, 0.272807018681266, 0.89317292560791095)
    59 = @(2.4647149518298099, 1.85934582936425, 0.187443549227999, 0.92393885590888403)
    61 = @(2.8395160547297298, 2.22124458039873, 0.12741283999594499, 0.91817141188113205)
    62 = @(2.9188234979498602, 2.1803938363077, 0.103056945310102, 0.90394116590684404)
    63 = @(2.5823256154854, 1.8989611182966899, 0.094213107779471894, 0.91816109273896196)
}

foreach ($row in $hpData.Keys) {
    $vals = $hpData[$row]
    $wsHP.Range("S$row").Value = $vals[0]
    $wsHP.Range("T$row").Value = $vals[1]
    $wsHP.Range("U$row").Value = $vals[2]
    $wsHP.Range("V$row").Value = $vals[3]
}

# -------------------------------------------------------------------------
# Sheet "CNNLSTM Temp": fill in previously-blank B:F raw-run columns for the
# S=3 / S=6 / S=12 blocks (the AVERAGE() in column G recalculates on its own).
# -------------------------------------------------------------------------
$wsCNN = $wb.Worksheets.Item("CNNLSTM Temp")

$cnnData = @{
    11 = @(1.01688854669771, 1.04141672564337, 1.5451864765862899, 0.92411909386145097, 1.23320673204392)
    12 = @(1.99352260367004, 1.60790309814746, 2.78227164021559, 2.00783101683637, 1.2105411746102399)
    13 = @(1.17875859590303, 1.4721097936174901, 2.4966902369772699, 1.31424542960843, 1.51475468628814)
    17 = @(0.908819494655417, 2.18714188215773, 2.7974925596742799, 2.7192225292466699, 1.62626626053167)
    18 = @(2.1687629736785401, 1.82792839812649, 3.0957506988729002, 4.1127081518374498, 2.6640810924228)
    19 = @(2.2153711570491001, 0.87115405008909597, 3.7649946979239899, 4.0296098879450399, 1.4424449661418)
    23 = @(2.59921667864832, 2.21835003678519, 4.3778570814684796, 2.6116290834220202, 2.3905273933246298)
    24 = @(3.1814293244382199, 2.64360050289649, 3.6721295150461, 2.9267447444583699, 2.1702134029101199)
    25 = @(2.4677236381093701, 2.1226413348920699, 2.98233479778066, 3.0943467900963402, 2.2445815165485299)
}

foreach ($row in $cnnData.Keys) {
    $vals = $cnnData[$row]
    $wsCNN.Range("B$row").Value = $vals[0]
    $wsCNN.Range("C$row").Value = $vals[1]
    $wsCNN.Range("D$row").Value = $vals[2]
    $wsCNN.Range("E$row").Value = $vals[3]
    $wsCNN.Range("F$row").Value = $vals[4]
}

# -------------------------------------------------------------------------
# Cosmetic view/selection state captured in the diff (best effort - window
# scroll position + active-cell selection on the two touched sheets).
# -------------------------------------------------------------------------
$wsHP.Activate()
$wsHP.Range("N49").Select()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 14
$wsHP.Range("S60").Select()

$wsCNN.Activate()
$wsCNN.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$wsCNN.Range("B26").Select()

$wsHP.Activate()
